# Updated solution for Tutorial 6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the date strings in column A (rows 3-21) from DD/MM/YYYY to DD-MM-YYYY.
# Force text interpretation (NumberFormat "@") while writing so that
# day-ambiguous values (e.g. 01-08-2022) aren't auto-parsed into real dates,
# then restore the default style so no stray formatting is left on the cell.
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    if ($current -ne $null) {
        $newValue = $current.Replace("/", "-")
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}

# Update attendance counts for specific rows (Real/Invalid/Absent recount)
$ws.Cells.Item(3, 4).Value = 1   # D3 Total Attendance Count
$ws.Cells.Item(3, 7).Value = 1   # G3 Invalid

$ws.Cells.Item(4, 4).Value = 1   # D4 Total Attendance Count
$ws.Cells.Item(4, 5).Value = 1   # E4 Real
$ws.Cells.Item(4, 8).Value = 0   # H4 Absent

$ws.Cells.Item(10, 4).Value = 1  # D10 Total Attendance Count
$ws.Cells.Item(10, 5).Value = 1  # E10 Real
$ws.Cells.Item(10, 8).Value = 0  # H10 Absent

$ws.Cells.Item(12, 4).Value = 1  # D12 Total Attendance Count
$ws.Cells.Item(12, 5).Value = 1  # E12 Real
$ws.Cells.Item(12, 8).Value = 0  # H12 Absent

$ws.Cells.Item(14, 4).Value = 1  # D14 Total Attendance Count
$ws.Cells.Item(14, 5).Value = 1  # E14 Real
$ws.Cells.Item(14, 8).Value = 0  # H14 Absent
